$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D35").Value = "Есть ЦА: Категория выстроена либо в единую линию, либо в две линии строго друг напротив друга (лицом друг к другу)"
$ws.Range("D37").Value = "Есть ЦА: Категория товаров для животных примыкает к ЦЕНТРАЛЬНОЙ АЛЛЕЕ и визуально доступна покупателям по ходу их движения без необходимости оборачиваться"
$ws.Range("D38").Value = "Есть ЦА: Категория товаров для животных примыкает к ПРОМО АЛЛЕЕ, находится дальше 5-ти метров от входа и визуально доступна покупателям по ходу их движения без необходимости оборачиваться"
